$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 785.2857
$ws.Range("I2").Value = 584.5
$ws.Range("J2").Value = 1053
$ws.Range("K2").Value = 584.5
$ws.Range("L2").Value = 1053
$ws.Range("M2").Value = -471.5
$ws.Range("N2").Value = -1279
$ws.Range("H6").Value = 113.333336
$ws.Range("I6").Value = 113.333336
$ws.Range("K6").Value = 340.000008
$ws.Range("M6").Value = -228.000008
$ws.Range("H17").Value = 984.6
$ws.Range("J17").Value = 892.1429000000001
$ws.Range("L17").Value = 2676.4287
$ws.Range("N17").Value = -3012.4287
$ws.Range("H38").Value = 694.3
$ws.Range("I38").Value = 241
$ws.Range("J38").Value = 2507.5
$ws.Range("K38").Value = 723
$ws.Range("L38").Value = 7522.5
$ws.Range("M38").Value = -351
$ws.Range("N38").Value = -8266.5
$ws.Range("H42").Value = 220.1
$ws.Range("I42").Value = 38.833332
$ws.Range("J42").Value = 492
$ws.Range("K42").Value = 116.499996
$ws.Range("L42").Value = 1476
$ws.Range("M42").Value = 113.500004
$ws.Range("N42").Value = -1936
$ws.Range("H58").Value = 2955.9092
$ws.Range("I58").Value = 728.75
$ws.Range("J58").Value = 4228.5713
$ws.Range("K58").Value = 2186.25
$ws.Range("L58").Value = 12685.7139
$ws.Range("M58").Value = -2036.25
$ws.Range("N58").Value = -12985.7139
$ws.Range("H86").Value = 214928.72
$ws.Range("I86").Value = 933
$ws.Range("J86").Value = 375425.5
$ws.Range("K86").Value = 933
$ws.Range("L86").Value = 375425.5
$ws.Range("M86").Value = 190
$ws.Range("N86").Value = -377671.5
$ws.Range("H89").Value = 214928.72
$ws.Range("I89").Value = 933
$ws.Range("J89").Value = 375425.5
$ws.Range("K89").Value = 4665
$ws.Range("L89").Value = 1877127.5
$ws.Range("M89").Value = 951
$ws.Range("N89").Value = -1888359.5
$ws.Range("H97").Value = 2898.4285
$ws.Range("J97").Value = 2898.4285
$ws.Range("L97").Value = 8695.2855
$ws.Range("N97").Value = -9687.2855
$ws.Range("H141").Value = 15000
$ws.Range("I141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("M141").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H13").Value = 7502625
$ws.Range("I13").Value = 15001000
$ws.Range("J13").Value = 4250
$ws.Range("K13").Value = 15001000
$ws.Range("L13").Value = 4250
$ws.Range("M13").Value = -15000856
$ws.Range("N13").Value = -4538
$ws.Range("H32").Value = 3457.5217
$ws.Range("I32").Value = 2162.8462
$ws.Range("J32").Value = 10670.714
$ws.Range("K32").Value = 2162.8462
$ws.Range("L32").Value = 10670.714
$ws.Range("M32").Value = -1875.8462
$ws.Range("N32").Value = -11244.714
$ws.Range("H45").Value = 2364.8667
$ws.Range("I45").Value = 1733.909
$ws.Range("K45").Value = 1733.909
$ws.Range("M45").Value = -1356.909
$ws.Range("H132").Value = 5382.7334
$ws.Range("J132").Value = 7598.8
$ws.Range("L132").Value = 22796.4
$ws.Range("N132").Value = -27856.4

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H11").Value = 1125
$ws.Range("I11").Value = 250
$ws.Range("K11").Value = 250
$ws.Range("M11").Value = -110
$ws.Range("H22").Value = 141911.11
$ws.Range("I22").Value = 147149.88
$ws.Range("J22").Value = 100001
$ws.Range("K22").Value = 147149.88
$ws.Range("L22").Value = 100001
$ws.Range("M22").Value = -146799.88
$ws.Range("N22").Value = -100701
$ws.Range("H86").Value = 4999
$ws.Range("I86").Value = 4999
$ws.Range("K86").Value = 4999
$ws.Range("M86").Value = -3876
$ws.Range("H89").Value = 4999
$ws.Range("I89").Value = 4999
$ws.Range("K89").Value = 24995
$ws.Range("M89").Value = -19379

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 3461286.8
$ws.Range("I4").Value = 4499539.5
$ws.Range("J4").Value = 444.66666
$ws.Range("K4").Value = 13498618.5
$ws.Range("L4").Value = 1333.99998
$ws.Range("M4").Value = -13498506.5
$ws.Range("N4").Value = -1557.99998
$ws.Range("H22").Value = 4830
$ws.Range("J22").Value = 4830
$ws.Range("L22").Value = 14490
$ws.Range("N22").Value = -14828
$ws.Range("H24").Value = 2061.5
$ws.Range("I24").Value = 175
$ws.Range("J24").Value = 2533.125
$ws.Range("K24").Value = 525
$ws.Range("L24").Value = 7599.375
$ws.Range("M24").Value = -295
$ws.Range("N24").Value = -8059.375
$ws.Range("H27").Value = 4830
$ws.Range("J27").Value = 4830
$ws.Range("L27").Value = 14490
$ws.Range("N27").Value = -14694
$ws.Range("H55").Value = 6883
$ws.Range("I55").Value = 2566.3333
$ws.Range("J55").Value = 8321.888999999999
$ws.Range("K55").Value = 7698.999899999999
$ws.Range("L55").Value = 24965.667
$ws.Range("M55").Value = -7521.999899999999
$ws.Range("N55").Value = -25319.667
$ws.Range("H81").Value = 299.8
$ws.Range("J81").Value = 300
$ws.Range("L81").Value = 900
$ws.Range("N81").Value = -3146
$ws.Range("H84").Value = 299.8
$ws.Range("J84").Value = 300
$ws.Range("L84").Value = 2700
$ws.Range("N84").Value = -13932
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()
$ws.Range("H126").Value = 1000
$ws.Range("I126").Value = 1000
$ws.Range("K126").Value = 3000
$ws.Range("M126").Value = 1940

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H25").Value = 1194.25
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 1194.25
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 1194.25
$ws.Range("M25").ClearContents()
$ws.Range("N25").Value = -2252.25
$ws.Range("H132").Value = 4488.0835
$ws.Range("I132").Value = 4641.727
$ws.Range("J132").Value = 2798
$ws.Range("K132").Value = 13925.181
$ws.Range("L132").Value = 8394
$ws.Range("M132").Value = -11395.181
$ws.Range("N132").Value = -13454

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 101373.4
$ws.Range("I82").Value = 1645
$ws.Range("K82").Value = 1645
$ws.Range("M82").Value = -1284
$ws.Range("H85").Value = 101373.4
$ws.Range("I85").Value = 1645
$ws.Range("K85").Value = 1645
$ws.Range("M85").Value = -397
$ws.Range("H93").Value = 866
$ws.Range("I93").Value = 866
$ws.Range("K93").Value = 866
$ws.Range("M93").Value = 382
$ws.Range("H122").Value = 6076.375
$ws.Range("I122").Value = 4890.4287
$ws.Range("J122").Value = 6998.778
$ws.Range("K122").Value = 14671.2861
$ws.Range("L122").Value = 20996.334
$ws.Range("M122").Value = -12221.2861
$ws.Range("N122").Value = -25896.334

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 1011250
$ws.Range("J26").Value = 1011250
$ws.Range("L26").Value = 1011250
$ws.Range("N26").Value = -1011836
$ws.Range("H122").Value = 2893.36
$ws.Range("I122").Value = 2580.125
$ws.Range("J122").Value = 3450.2222
$ws.Range("K122").Value = 7740.375
$ws.Range("L122").Value = 10350.6666
$ws.Range("M122").Value = -5290.375
$ws.Range("N122").Value = -15250.6666
$ws.Range("H126").Value = 1400.2
$ws.Range("I126").Value = 1222.4445
$ws.Range("K126").Value = 3667.3335
$ws.Range("M126").Value = -1197.3335
